$wb = $excel.ActiveWorkbook

# --- Sheet1: "Probability Cost" data updates ---
$ws1 = $wb.Worksheets.Item("Probability Cost")

# Row 7 (piece count = 6): B7 0.25 -> 0.3 ; D7 0.3 -> 0.25
$ws1.Range("B7").Value = 0.3
$ws1.Range("D7").Value = 0.25

# Row 8 (piece count = 7): C8 0.3 -> 0.35 ; D8 0.4 -> 0.35
$ws1.Range("C8").Value = 0.35
$ws1.Range("D8").Value = 0.35

# Row 10 (piece count = 9): C10 0.15 -> 0.2 ; D10 0.2 -> 0.25
$ws1.Range("C10").Value = 0.2
$ws1.Range("D10").Value = 0.25

# --- Sheet2: "Max N" worksheet reference (no content changes besides view) ---
$ws2 = $wb.Worksheets.Item("Max N")

# --- Sheet view / selection updates ---
$ws1.Activate()
$ws1.Range("G10").Select()

$ws2.Activate()
$ws2.Range("F2").Select()

$ws1.Activate()
